$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Fix the missing accent: "correspondent a notre" -> "correspondent à notre"
#    (done first, while the sentence is still a single run, so the
#    paragraph's run-list gets rebuilt with the corrected text before
#    we start slicing it into separate runs below)
# ------------------------------------------------------------------
$d.Content.Find.Execute("correspondent a notre", $true, $false, $false, $false, $false, $true, 1, $false, "correspondent à notre", 2)

# ------------------------------------------------------------------
# 2. Split "...Git. Sélectionne" | "r les champs..." into two runs.
#    Adding (and immediately deleting) a temporary bookmark at the
#    desired split point forces the host to break the run there
#    without altering any text.
# ------------------------------------------------------------------
$rngSplit1 = $d.Content
$rngSplit1.Find.Execute("Git. Sélectionne", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngSplit1.Collapse(0)
$d.Bookmarks.Add("TempSplit1", $rngSplit1)
$d.Bookmarks("TempSplit1").Delete()

# ------------------------------------------------------------------
# 3. Split "Derniere étape" | " pour la création..." into two runs,
#    the same way.
# ------------------------------------------------------------------
$rngSplit2 = $d.Content
$rngSplit2.Find.Execute("Derniere étape", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngSplit2.Collapse(0)
$d.Bookmarks.Add("TempSplit2", $rngSplit2)
$d.Bookmarks("TempSplit2").Delete()

# ------------------------------------------------------------------
# 4. Move the "_GoBack" bookmark (Word's "last edit" marker) to right
#    after "...correspondent à", which is where it ends up once the
#    text above has been edited. Re-adding a bookmark with a name
#    that already exists elsewhere in the document relocates it, so
#    this also removes it from its old spot near the end of the doc.
# ------------------------------------------------------------------
$rngGoBack = $d.Content
$rngGoBack.Find.Execute("correspondent à", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngGoBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngGoBack)
